# Testersky slovnik pokracovanie - add "Stacktrace" glossary entry.
#
# The document ends with two empty paragraphs right after the "--" separator
# paragraph. The first of those two empty paragraphs gets replaced with a new
# glossary entry consisting of two paragraphs:
#   1) a bold "Stacktrace " heading (with spell-check markers around the
#      single word "Stacktrace", matching the style used elsewhere in the
#      glossary)
#   2) the definition paragraph, with the phrase "retazec volani funkcii"
#      rendered in bold in the middle of the sentence.
# The final empty paragraph (end of document / before sectPr) is left as-is.

$d = $word.ActiveDocument

# The document ends with two empty paragraphs (the very last paragraph of the
# story, and the one right before it). We want to replace the *first* of
# those two empty trailing paragraphs -- i.e. the second-to-last paragraph
# overall -- which sits right after the "--" separator paragraph.
$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n - 1)

if ($target.Range.Text -ne "`r") {
    throw "Expected the second-to-last paragraph to be empty, found: '$($target.Range.Text)'"
}

$r = $target.Range

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $w + '>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Stacktrace</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '</w:p>'

$para2 = '<w:p ' + $w + '>' +
         '<w:r><w:t xml:space="preserve">(alebo „výpis zásobníka volaní“) je technický výpis, ktorý ukazuje </w:t></w:r>' +
         '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>reťazec volaní funkcií</w:t></w:r>' +
         '<w:r><w:t>, ktoré viedli k chybe alebo výnimke v programe.</w:t></w:r>' +
         '</w:p>'

$body = $para1 + $para2

$flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document ' + $w + '><w:body>' + $body + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($flatOpc)

Write-Host "Inserted Stacktrace glossary entry."
